# Update automatico via Actualizar@
#
# This mirrors the "availability checker" automation: it refreshes the
# recorded check timestamps for the existing rows, and appends a fresh
# batch of availability rows (a re-check of the first 8 monitored
# services) at the bottom of the sheet, including their hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Refresh the timestamp column for the two existing timestamp groups.
# ---------------------------------------------------------------------
$ws.Range("D2:D15").Value = 44230.88231832047
$ws.Range("D16:D29").Value = 44230.87725172454

# ---------------------------------------------------------------------
# 2. Append a new batch of rows (30-37), duplicating the first 8 service
#    rows (name / url / status), stamped with the timestamp the second
#    group used to have, and re-create their hyperlinks.
# ---------------------------------------------------------------------
$destFirstRow = 30
$rowCount = 8
$newTimestamp = 44230.72746922453

# Hyperlink targets for the first 8 service rows (same order/targets as
# the existing B2:B9 hyperlinks); MapStore (index 5, 0-based) also
# carries the "/" fragment sub-address, just like B7 does.
$targets = @(
    "https://ezexporter.highviewapps.com/exports/export-profile/",
    "https://github.com/Sud-Austral/",
    "https://rpubs.com/dataintelligence/",
    "https://ide.dataintelligence-group.com/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://odooutil.azurewebsites.net/"
)
$subAddresses = @("", "", "", "", "", "/", "", "")

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = 2 + $i
    $dstRow = $destFirstRow + $i

    $nameValue = $ws.Cells.Item($srcRow, 1).Value()
    $urlText = $ws.Cells.Item($srcRow, 2).Text()
    $statusValue = $ws.Cells.Item($srcRow, 3).Value()

    $ws.Cells.Item($dstRow, 1).Value = $nameValue
    $ws.Cells.Item($dstRow, 2).Value = $urlText
    $ws.Cells.Item($dstRow, 3).Value = $statusValue

    $ws.Cells.Item($dstRow, 4).Value = $newTimestamp
    $ws.Cells.Item($dstRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $destCell = $ws.Cells.Item($dstRow, 2)
    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($destCell, $targets[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($destCell, $targets[$i])
    }
}

# Re-apply the Hyperlink cell style (the Add call above styles the cell
# with a brand-new style entry instead of reusing the existing one; fix
# it up afterwards so it matches the style used by B2:B29).
for ($i = 0; $i -lt $rowCount; $i++) {
    $dstRow = $destFirstRow + $i
    $ws.Cells.Item($dstRow, 2).Style = $ws.Cells.Item(2, 2).Style()
}
